$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-22 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-23 Friday", 2) | Out-Null
$d.Content.Find.Execute("56×54=3024", $true, $false, $false, $false, $false, $true, 1, $false, "80×29=2320", 2) | Out-Null
$d.Content.Find.Execute("32×71=2272", $true, $false, $false, $false, $false, $true, 1, $false, "25×22=550", 2) | Out-Null
$d.Content.Find.Execute("43×23=989", $true, $false, $false, $false, $false, $true, 1, $false, "89×97=8633", 2) | Out-Null
$d.Content.Find.Execute("84×52=4368", $true, $false, $false, $false, $false, $true, 1, $false, "41×43=1763", 2) | Out-Null
$d.Content.Find.Execute("99×60=5940", $true, $false, $false, $false, $false, $true, 1, $false, "88×94=8272", 2) | Out-Null
$d.Content.Find.Execute("82×96=7872", $true, $false, $false, $false, $false, $true, 1, $false, "12×49=588", 2) | Out-Null
$d.Content.Find.Execute("67×26=1742", $true, $false, $false, $false, $false, $true, 1, $false, "45×66=2970", 2) | Out-Null
$d.Content.Find.Execute("17×15=255", $true, $false, $false, $false, $false, $true, 1, $false, "43×57=2451", 2) | Out-Null
$d.Content.Find.Execute("13×13=169", $true, $false, $false, $false, $false, $true, 1, $false, "39×52=2028", 2) | Out-Null
$d.Content.Find.Execute("86×26=2236", $true, $false, $false, $false, $false, $true, 1, $false, "98×71=6958", 2) | Out-Null
$d.Content.Find.Execute("93×29=2697", $true, $false, $false, $false, $false, $true, 1, $false, "40×35=1400", 2) | Out-Null
$d.Content.Find.Execute("26×84=2184", $true, $false, $false, $false, $false, $true, 1, $false, "69×68=4692", 2) | Out-Null
$d.Content.Find.Execute("40×87=3480", $true, $false, $false, $false, $false, $true, 1, $false, "65×98=6370", 2) | Out-Null
$d.Content.Find.Execute("43×30=1290", $true, $false, $false, $false, $false, $true, 1, $false, "12×42=504", 2) | Out-Null
$d.Content.Find.Execute("50×53=2650", $true, $false, $false, $false, $false, $true, 1, $false, "11×30=330", 2) | Out-Null
$d.Content.Find.Execute("64×41=2624", $true, $false, $false, $false, $false, $true, 1, $false, "40×50=2000", 2) | Out-Null
$d.Content.Find.Execute("74×58=4292", $true, $false, $false, $false, $false, $true, 1, $false, "20×47=940", 2) | Out-Null
$d.Content.Find.Execute("17×35=595", $true, $false, $false, $false, $false, $true, 1, $false, "50×42=2100", 2) | Out-Null
$d.Content.Find.Execute("97×25=2425", $true, $false, $false, $false, $false, $true, 1, $false, "65×84=5460", 2) | Out-Null
$d.Content.Find.Execute("69×26=1794", $true, $false, $false, $false, $false, $true, 1, $false, "77×14=1078", 2) | Out-Null
$d.Content.Find.Execute("67×80=5360", $true, $false, $false, $false, $false, $true, 1, $false, "24×47=1128", 2) | Out-Null
$d.Content.Find.Execute("69×80=5520", $true, $false, $false, $false, $false, $true, 1, $false, "75×63=4725", 2) | Out-Null
$d.Content.Find.Execute("26×43=1118", $true, $false, $false, $false, $false, $true, 1, $false, "56×49=2744", 2) | Out-Null
$d.Content.Find.Execute("36×57=2052", $true, $false, $false, $false, $false, $true, 1, $false, "56×17=952", 2) | Out-Null
$d.Content.Find.Execute("16×22=352", $true, $false, $false, $false, $false, $true, 1, $false, "36×30=1080", 2) | Out-Null
